# Fix the car title for row 3 on "findNewCarTest": "Ho Cars" -> "Honda Cars"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("findNewCarTest")
$ws1.Range("D3").Value = "Honda Cars"

# Make "findNewCarTest" the active sheet/tab and move its selection to H6
$ws1.Activate()
$ws1.Range("H6").Select()
